$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "34.116.04"
Set-TextValue "E2" "  +0.13%  "
Set-TextValue "D3" "1.790.58"
Set-TextValue "E3" "  -0.11%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "227.16"
Set-TextValue "E5" "  +1.41%  "
Set-TextValue "E6" "  -0.75%  "
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "32.38"
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "E9" "  +4.39%  "
Set-TextValue "E10" "  -2.77%  "
Set-TextValue "D11" "0.0942"
Set-TextValue "E11" "  +1.30%  "
Set-TextValue "D12" "2.049.58"
Set-TextValue "E12" "  +0.00%  "
Set-TextValue "D13" "11.50"
Set-TextValue "E13" "  +5.97%  "
Set-TextValue "D14" "1.798.77"
Set-TextValue "E14" "  +0.33%  "
Set-TextValue "D15" "0.623"
Set-TextValue "E15" "  -0.12%  "
Set-TextValue "D16" "34.114.43"
Set-TextValue "E16" "  +0.17%  "
Set-TextValue "D17" "4.18"
Set-TextValue "E17" "  +0.55%  "
Set-TextValue "D18" "68.00"
Set-TextValue "E18" "  +0.10%  "
Set-TextValue "D19" "243.72"
Set-TextValue "E19" "  +0.19%  "
Set-TextValue "E20" "  -0.57%  "
Set-TextValue "D21" "1.00"
Set-TextValue "E21" "  -0.13%  "
Set-TextValue "D22" "10.88"
Set-TextValue "D23" "4.11"
Set-TextValue "E23" "  +1.00%  "
Set-TextValue "D24" "2.06"
Set-TextValue "E24" "  -2.11%  "
Set-TextValue "D25" "162.03"
Set-TextValue "E25" "  +2.04%  "
Set-TextValue "E26" "  +2.49%  "
Set-TextValue "D27" "16.26"
Set-TextValue "E27" "  -0.11%  "
Set-TextValue "E28" "  +1.34%  "
Set-TextValue "E29" "  +0.13%  "
Set-TextValue "E30" "  +2.53%  "
Set-TextValue "D31" "0.0517"
Set-TextValue "E31" "  +0.13%  "
Set-TextValue "D32" "3.66"
Set-TextValue "E32" "  +0.38%  "
Set-TextValue "D33" "3.63"
Set-TextValue "E33" "  +4.19%  "
Set-TextValue "E34" "  +1.79%  "
Set-TextValue "D35" "1.407.13"
Set-TextValue "E35" "  +1.57%  "
Set-TextValue "E36" "  +1.57%  "
Set-TextValue "E37" "  -0.02%  "
Set-TextValue "D38" "0.0189"
Set-TextValue "E38" "  +2.75%  "
Set-TextValue "E39" "  +8.83%  "
Set-TextValue "D40" "80.12"
Set-TextValue "E40" "  +1.22%  "
Set-TextValue "E41" "  +0.56%  "
Set-TextValue "D42" "0.923"
Set-TextValue "E42" "  +1.05%  "
Set-TextValue "D43" "2.71"
Set-TextValue "E43" "  +0.10%  "
Set-TextValue "D44" "13.36"
Set-TextValue "E44" "  +12.03%  "
Set-TextValue "D45" "0.0₆0142"
Set-TextValue "E45" "  +4.32%  "
Set-TextValue "D46" "6.09"
Set-TextValue "E46" "  +4.38%  "
Set-TextValue "B47" "WEMIXToken"
Set-TextValue "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D47" "1.08"
Set-TextValue "E47" "  +2.47%  "
Set-TextValue "B48" "Kaspa"
Set-TextValue "C48" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D48" "0.0506"
Set-TextValue "E48" "  +1.89%  "
Set-TextValue "D49" "107.48"
Set-TextValue "E49" "  +0.27%  "
Set-TextValue "D50" "1.951.77"
Set-TextValue "E50" "  +0.20%  "
Set-TextValue "E51" "  +0.04%  "
